# Update cached market-price / profit figures on several leve sheets.
# Each entry: (SheetName, Row, Column, NewValue)
$changes = @(
    @("ALC", 58, 8, 21812.549),
    @("ALC", 58, 9, 320),
    @("ALC", 58, 10, 24148.695),
    @("ALC", 58, 11, 960),
    @("ALC", 58, 12, 72446.08499999999),
    @("ALC", 58, 13, -810),
    @("ALC", 58, 14, -72746.08499999999),
    @("ALC", 69, 8, 3554.75),
    @("ALC", 69, 10, 3570.3333),
    @("ALC", 69, 12, 10710.9999),
    @("ALC", 69, 14, -12458.9999),
    @("ALC", 72, 8, 3554.75),
    @("ALC", 72, 10, 3570.3333),
    @("ALC", 72, 12, 32132.9997),
    @("ALC", 72, 14, -40868.9997),
    @("ALC", 82, 8, 5388.2),
    @("ALC", 82, 9, 3647),
    @("ALC", 82, 11, 10941),
    @("ALC", 82, 13, -10535),
    @("ALC", 85, 8, 5388.2),
    @("ALC", 85, 9, 3647),
    @("ALC", 85, 11, 10941),
    @("ALC", 85, 13, -9537),
    @("ALC", 135, 8, 899.625),
    @("ALC", 135, 9, 799.2308),
    @("ALC", 135, 11, 7193.077200000001),
    @("ALC", 135, 13, -4658.077200000001),
    @("BSM", 105, 8, 1838.4166),
    @("BSM", 105, 9, 1758.1333),
    @("BSM", 105, 11, 1758.1333),
    @("BSM", 105, 13, -11.13329999999996),
    @("BSM", 134, 8, 1439.0182),
    @("BSM", 134, 9, 992.5238000000001),
    @("BSM", 134, 10, 2881.5386),
    @("BSM", 134, 11, 2977.5714),
    @("BSM", 134, 12, 8644.6158),
    @("BSM", 134, 13, -442.5714000000003),
    @("BSM", 134, 14, -13714.6158),
    @("CRP", 31, 8, 5114.8276),
    @("CRP", 31, 9, 3612.3333),
    @("CRP", 31, 11, 3612.3333),
    @("CRP", 31, 13, -3317.3333),
    @("CRP", 34, 8, 5114.8276),
    @("CRP", 34, 9, 3612.3333),
    @("CRP", 34, 11, 3612.3333),
    @("CRP", 34, 13, -3410.3333),
    @("CRP", 58, 8, 11907626),
    @("CRP", 58, 9, 1643.5862),
    @("CRP", 58, 10, 38467130),
    @("CRP", 58, 11, 1643.5862),
    @("CRP", 58, 12, 38467130),
    @("CRP", 58, 13, -1440.5862),
    @("CRP", 58, 14, -38467536),
    @("CRP", 122, 8, 3567.318),
    @("CRP", 122, 9, 4158),
    @("CRP", 122, 10, 3393.5881),
    @("CRP", 122, 11, 12474),
    @("CRP", 122, 12, 10180.7643),
    @("CRP", 122, 13, -10024),
    @("CRP", 122, 14, -15080.7643),
    @("CRP", 136, 8, 11907626),
    @("CRP", 136, 9, 1643.5862),
    @("CRP", 136, 10, 38467130),
    @("CRP", 136, 11, 4930.7586),
    @("CRP", 136, 12, 115401390),
    @("CRP", 136, 13, -2380.7586),
    @("CRP", 136, 14, -115406490),
    @("CUL", 3, 8, 4020.5454),
    @("CUL", 3, 9, 2272),
    @("CUL", 3, 10, 5477.6665),
    @("CUL", 3, 11, 6816),
    @("CUL", 3, 12, 16432.9995),
    @("CUL", 3, 13, -6704),
    @("CUL", 3, 14, -16656.9995),
    @("CUL", 75, 8, 2332.8),
    @("CUL", 75, 9, 1012.6667),
    @("CUL", 75, 10, 2898.5715),
    @("CUL", 75, 11, 3038.0001),
    @("CUL", 75, 12, 8695.7145),
    @("CUL", 75, 13, -2040.0001),
    @("CUL", 75, 14, -10691.7145),
    @("CUL", 78, 8, 2332.8),
    @("CUL", 78, 9, 1012.6667),
    @("CUL", 78, 10, 2898.5715),
    @("CUL", 78, 11, 9114.0003),
    @("CUL", 78, 12, 26087.1435),
    @("CUL", 78, 13, -4122.0003),
    @("CUL", 78, 14, -36071.1435),
    @("CUL", 87, 8, 11700),
    @("CUL", 87, 9, 6660),
    @("CUL", 87, 11, 19980),
    @("CUL", 87, 13, -18732),
    @("CUL", 90, 8, 11700),
    @("CUL", 90, 9, 6660),
    @("CUL", 90, 11, 59940),
    @("CUL", 90, 13, -53700),
    @("CUL", 120, 8, 19290),
    @("CUL", 120, 9, 18757.5),
    @("CUL", 120, 11, 56272.5),
    @("CUL", 120, 13, -51434.5),
    @("CUL", 125, 8, 2406),
    @("CUL", 125, 9, 1515),
    @("CUL", 125, 10, 3000),
    @("CUL", 125, 11, 4545),
    @("CUL", 125, 12, 9000),
    @("CUL", 125, 13, 375),
    @("CUL", 125, 14, -18840),
    @("CUL", 126, 8, 1638.3334),
    @("CUL", 126, 9, 943.3333),
    @("CUL", 126, 10, 2333.3333),
    @("CUL", 126, 11, 2829.9999),
    @("CUL", 126, 12, 6999.999899999999),
    @("CUL", 126, 13, 2110.0001),
    @("CUL", 126, 14, -16879.9999),
    @("CUL", 130, 8, 2062.5),
    @("CUL", 130, 9, 750),
    @("CUL", 130, 10, 2500),
    @("CUL", 130, 11, 2250),
    @("CUL", 130, 12, 7500),
    @("CUL", 130, 13, 2770),
    @("CUL", 130, 14, -17540),
    @("CUL", 133, 8, 6231.778),
    @("CUL", 133, 9, 9055),
    @("CUL", 133, 10, 3973.2),
    @("CUL", 133, 11, 27165),
    @("CUL", 133, 12, 11919.6),
    @("CUL", 133, 13, -22105),
    @("CUL", 133, 14, -22039.6),
    @("CUL", 134, 8, 2371.5908),
    @("CUL", 134, 9, 1324.75),
    @("CUL", 134, 10, 2969.7856),
    @("CUL", 134, 11, 3974.25),
    @("CUL", 134, 12, 8909.356800000001),
    @("CUL", 134, 13, 1095.75),
    @("CUL", 134, 14, -19049.3568),
    @("CUL", 136, 8, 2285.1667),
    @("CUL", 136, 9, 1453.0625),
    @("CUL", 136, 10, 3949.375),
    @("CUL", 136, 11, 4359.1875),
    @("CUL", 136, 12, 11848.125),
    @("CUL", 136, 13, 740.8125),
    @("CUL", 136, 14, -22048.125),
    @("CUL", 137, 8, 2051.1292),
    @("CUL", 137, 9, 1454.0625),
    @("CUL", 137, 10, 2688),
    @("CUL", 137, 11, 4362.1875),
    @("CUL", 137, 12, 8064),
    @("CUL", 137, 13, 737.8125),
    @("CUL", 137, 14, -18264),
    @("CUL", 138, 8, 2249.0715),
    @("CUL", 138, 9, 1238.7142),
    @("CUL", 138, 10, 3259.4285),
    @("CUL", 138, 11, 3716.1426),
    @("CUL", 138, 12, 9778.2855),
    @("CUL", 138, 13, 1423.8574),
    @("CUL", 138, 14, -20058.2855),
    @("CUL", 139, 8, 6415215.5),
    @("CUL", 139, 9, 9261419),
    @("CUL", 139, 10, 11257.417),
    @("CUL", 139, 11, 27784257),
    @("CUL", 139, 12, 33772.251),
    @("CUL", 139, 13, -27779117),
    @("CUL", 139, 14, -44052.251),
    @("CUL", 140, 8, 6948073),
    @("CUL", 140, 9, 41667676),
    @("CUL", 140, 10, 4152.5),
    @("CUL", 140, 11, 125003028),
    @("CUL", 140, 12, 12457.5),
    @("CUL", 140, 13, -124997848),
    @("CUL", 140, 14, -22817.5),
    @("GSM", 122, 8, 3708.0645),
    @("GSM", 122, 9, 2725),
    @("GSM", 122, 10, 5069.231),
    @("GSM", 122, 11, 8175),
    @("GSM", 122, 12, 15207.693),
    @("GSM", 122, 13, -5725),
    @("GSM", 122, 14, -20107.693),
    @("LTW", 133, 8, 29642.715),
    @("LTW", 133, 10, 29642.715),
    @("LTW", 133, 12, 29642.715),
    @("LTW", 133, 14, -34702.715),
    @("LTW", 135, 8, 29597.785),
    @("LTW", 135, 10, 29597.785),
    @("LTW", 135, 12, 29597.785),
    @("LTW", 135, 14, -39737.785),
    @("LTW", 136, 8, 1664.9231),
    @("LTW", 136, 9, 1151.1428),
    @("LTW", 136, 11, 3453.4284),
    @("LTW", 136, 13, -903.4284000000002),
    @("WVR", 46, 8, 37771.285),
    @("WVR", 46, 10, 37771.285),
    @("WVR", 46, 12, 37771.285),
    @("WVR", 46, 14, -38233.285),
    @("WVR", 81, 8, 2780.8572),
    @("WVR", 81, 9, 1577.3334),
    @("WVR", 81, 10, 10002),
    @("WVR", 81, 11, 3154.6668),
    @("WVR", 81, 12, 20004),
    @("WVR", 81, 13, -2093.6668),
    @("WVR", 81, 14, -22126),
    @("WVR", 84, 8, 2780.8572),
    @("WVR", 84, 9, 1577.3334),
    @("WVR", 84, 10, 10002),
    @("WVR", 84, 11, 15773.334),
    @("WVR", 84, 12, 100020),
    @("WVR", 84, 13, -10469.334),
    @("WVR", 84, 14, -110628),
    @("WVR", 134, 8, 37771.285),
    @("WVR", 134, 10, 37771.285),
    @("WVR", 134, 12, 113313.855),
    @("WVR", 134, 14, -118383.855),
    @("WVR", 136, 8, 1219.9667),
    @("WVR", 136, 9, 727.2273),
    @("WVR", 136, 11, 2181.6819),
    @("WVR", 136, 13, 368.3181)
)

$wb = $excel.ActiveWorkbook

$wsCache = @{}

foreach ($change in $changes) {
    $sheetName = $change[0]
    $row = [int]$change[1]
    $col = [int]$change[2]
    $val = [double]$change[3]

    if (-not $wsCache.ContainsKey($sheetName)) {
        $wsCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $wsCache[$sheetName]

    $ws.Cells.Item($row, $col).Value = $val
}
